$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Jagged array holding the full target values for range A1:Q25
$jagged = @(
        @($null, 0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15),
        @(0, 5.830559078307802, 0, 4.05224556568083, 14.65450099282531, 24.98979208541448, 35.16368355843755, 1.877117595328463, 3.597007928081834, 11.16073994883218, 15.81134233393687, 0, 16.77717827336655, 0, 16.21839987150576, 0, 0),
        @(1, 5.529312788888233, 0, 4.022574270731143, 13.83323472501031, 24.0353831458494, 33.64718806081251, 1.626852739916426, 3.721424056453005, 10.97869074408844, 15.92954408615093, 0, 15.7854127867996, 0, 15.2364655305274, 0, 0),
        @(2, 5.331220256968415, 0, 4.003554407556764, 13.30304726557518, 23.44176032945489, 32.69707681081033, 1.655212418641725, 3.801002351107744, 10.86959112477356, 16.00180997569337, 0, 15.14793366132119, 0, 14.60392844432243, 0, 0),
        @(3, 5.237867756842999, 0, 3.996204728907074, 13.07657815407962, 23.18347977115934, 32.27890340071477, 1.722856744490209, 3.836969538647418, 10.82059244392093, 16.02064563772023, 0, 14.88889913189589, 0, 14.34288476103972, 0, 0),
        @(4, 5.210713345490284, 0, 3.995717432942752, 13.03378267761558, 23.12244628458413, 32.17656764359949, 1.734264673933993, 3.846466446947309, 10.80611215159094, 16.01093843705722, 0, 14.85508690716737, 0, 14.30410393181637, 0, 0),
        @(5, 5.29918323802446, 0, 4.0055276449502, 13.28718156858112, 23.38966235169254, 32.60369070053432, 1.656371996951416, 3.810667046816126, 10.85161103530296, 15.96701281585266, 0, 15.17023312534914, 0, 14.61380353641702, 0, 0),
        @(6, 5.690979778880837, 0, 4.04509761452183, 14.36089661751481, 24.60134872934431, 34.53533061713146, 1.791554200614444, 3.650814437017921, 11.07510292837192, 15.80601888577462, 0, 16.47298325880497, 0, 15.90239042522259, 0, 0),
        @(7, 6.404317403571902, 0, 4.112298030398885, 16.28096531441616, 26.98659040516356, 38.30040908245903, 2.389981713538565, 3.347158720456562, 11.56422721109525, 15.54980085561083, 0, 18.7584884177765, 0, 18.1707973006373, 0, 0),
        @(8, 6.894766808096158, 0, 4.180328325712438, 17.04964720568243, 28.416357772338, 40.47280718988566, 2.783468161794678, 3.147757274602025, 11.83531140154769, 15.24523832039933, 0, 20.32118689160265, 0, 19.41997884896279, 0, 0),
        @(9, 7.45383994630948, 0, 4.439035316276767, 13.2506726045486, 26.87996532716276, 37.57587436205231, 3.420531357889032, 3.108534114314186, 11.13624990958127, 14.2947863011101, 0, 21.14074018419854, 0, 17.60641631982411, 0, 0),
        @(10, 7.87224031267665, 0, 4.682328454435482, 10.49749067168011, 25.32156709228715, 34.74703321838306, 4.501749428747595, 3.096234691884984, 10.5064490852419, 13.71770176520154, 0, 21.489155413723, 0, 15.85956532270564, 0, 0),
        @(11, 8.183296071580163, 0, 4.923888926667472, 8.847219419433397, 23.54593350274641, 31.60372715831593, 5.719049034390846, 3.118385818513707, 9.861047841945709, 13.35298182083137, 0, 21.53726905560676, 0, 14.01698511566279, 0, 0),
        @(12, 8.351983236074771, 0, 5.095342434883233, 8.703445370665877, 22.17436313665178, 29.20287044102709, 6.61854702327035, 3.151888873890152, 9.398534787502536, 13.21190235641681, 0, 21.44297428299906, 0, 12.66964872327842, 0, 0),
        @(13, 8.357228487507083, 0, 5.133607564586613, 8.803357614647558, 21.77740119885651, 28.51971912830982, 6.823574400113162, 3.172202923281451, 9.278234922173905, 13.20494740108208, 0, 21.36294000771703, 0, 12.31286698061128, 0, 0),
        @(14, 8.087178739290483, 0, 5.053188088114513, 8.604149807735487, 21.39465439728346, 27.97274897551743, 6.545362047991158, 3.257936214250432, 9.269654096663849, 13.39878772095444, 0, 20.73175930957046, 0, 12.02898201959384, 0, 0),
        @(15, 7.792474620350902, 0, 4.896852096572728, 8.326804143961333, 21.85917549495649, 28.88664399322983, 5.802282015481826, 3.306969062670058, 9.510830113549172, 13.60807472593944, 0, 20.2901382586493, 0, 12.54945744869561, 0, 0),
        @(16, 7.461749104010942, 0, 4.667898617079578, 8.978674299015038, 23.13427700329596, 31.21463740710906, 4.623283059843406, 3.317921701395125, 10.00748450855384, 13.93343837490808, 0, 19.95343772285617, 0, 13.85771459859185, 0, 0),
        @(17, 7.084181101253872, 0, 4.432287621882044, 11.41159229027049, 24.87128387600407, 34.33502729851654, 3.367405994747044, 3.312859456203136, 10.64872961890068, 14.39120708228806, 0, 19.7756628834712, 0, 15.75968341740445, 0, 0),
        @(18, 6.694136191424567, 0, 4.171582858995226, 16.79678757028747, 27.90689935889245, 39.66551864967138, 2.678681277476407, 3.229060351002395, 11.71030533603315, 15.21612493661112, 0, 19.98479975991595, 0, 19.11913979642922, 0, 0),
        @(19, 6.999874148678837, 0, 4.187962648136119, 18.15340549756189, 29.39878096902437, 42.03082275968226, 3.014132311406756, 3.068308916366282, 12.07808302731794, 15.15370091793979, 0, 21.10451589156287, 0, 20.45436066473625, 0, 0),
        @(20, 7.240690793122894, 0, 4.205540783741238, 18.7650878375077, 30.28233261135528, 43.41035229686437, 3.216873582823671, 2.955437062254815, 12.29226213871, 15.1079758097338, 0, 21.79367494020267, 0, 21.1452219160948, 0, 0),
        @(21, 7.141453691681112, 0, 4.193183438910901, 18.45191808456669, 29.85692143135218, 42.75571556182251, 3.10911786130864, 3.003806814828043, 12.1955550043558, 15.17264702114292, 0, 21.40535786712981, 0, 20.76740670452914, 0, 0),
        @(22, 6.698463739141316, 0, 4.149689287236222, 17.20421086847962, 28.15237299728983, 40.10675752087091, 2.694306009016411, 3.210312677989294, 11.80661126316915, 15.35432479548331, 0, 19.91929088799989, 0, 19.30086016754908, 0, 0),
        @(23, 6.174332076697759, 0, 4.099103963959359, 15.76391874505065, 26.27293278441696, 37.16385096739651, 2.231658186887858, 3.443031727872239, 11.40042036510001, 15.55719741460978, 0, 18.20432403499387, 0, 17.60296405012826, 0, 0)
)

$rows = $jagged.Count
$cols = $jagged[0].Count
$data = New-Object 'object[,]' $rows,$cols
for ($r = 0; $r -lt $rows; $r++) {
    for ($c = 0; $c -lt $cols; $c++) {
        $data[$r,$c] = $jagged[$r][$c]
    }
}

$ws.Range("A1:Q25").Value = $data

# Apply the header style (bold, bordered, centered) used for the rest of the
# header row to the two newly added header cells (P1, Q1) without disturbing
# the values that were just written.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
